$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 69; Excel shifts rows 69:141 down to 70:142
# and copies formatting (e.g. the date style on column D) from the row
# being pushed down.
$ws.Rows("69:69").Insert()

# Populate the newly inserted row 69 with the new weekly record.
# Columns A, B, C, E, F, G, H, I, R keep the same values as the
# (now shifted) row below, so only set the columns that differ.
$ws.Range("A69").Value = 7
$ws.Range("B69").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C69").Value = "Ñuble"
$ws.Range("D69").Value = 45068
$ws.Range("E69").Value = 16
$ws.Range("F69").Value = 100112037
$ws.Range("G69").Value = "Cebollín"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 200
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 7000
$ws.Range("M69").Value = 6500
$ws.Range("N69").Value = "$/paquete 36 unidades"
$ws.Range("O69").Value = "Provincia de Diguillín"
$ws.Range("P69").Value = 181
$ws.Range("Q69").Value = 36
$ws.Range("R69").Value = "Hortaliza"
